$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.130
Model:                            OLS   Adj. R-squared:                  0.091
Method:                 Least Squares   F-statistic:                     3.343
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           4.01e-07
Time:                        16:11:18   Log-Likelihood:                -303.41
No. Observations:                 537   AIC:                             654.8
Df Residuals:                     513   BIC:                             757.7
Df Model:                          23                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.1282      0.375      3.011      0.003       0.392       1.864
NEMScr               -0.0022      0.001     -1.903      0.058      -0.005    7.23e-05
Ranking               0.0013      0.001      1.399      0.162      -0.001       0.003
LangScr               0.0010      0.000      2.790      0.005       0.000       0.002
MathScr              -0.0021      0.001     -3.976      0.000      -0.003      -0.001
ScienScr          -6.391e-05      0.000     -0.155      0.877      -0.001       0.001
SchoolRegion_1        0.0927      0.169      0.549      0.583      -0.239       0.424
SchoolRegion_2        0.2705      0.155      1.744      0.082      -0.034       0.575
SchoolRegion_3        0.1979      0.188      1.056      0.292      -0.170       0.566
EdTypeCode_1         -0.0624      0.111     -0.564      0.573      -0.280       0.155
EdTypeCode_2          0.1246      0.106      1.176      0.240      -0.084       0.333
SchoolType_1         -0.0099      0.112     -0.088      0.930      -0.229       0.209
SchoolType_2         -0.1418      0.096     -1.474      0.141      -0.331       0.047
SchoolType_3         -0.0679      0.091     -0.743      0.458      -0.247       0.111
SchoolType_4         -0.2123      0.282     -0.753      0.452      -0.766       0.342
MotherEd_1            0.0757      0.102      0.743      0.458      -0.125       0.276
MotherEd_2            0.0925      0.088      1.048      0.295      -0.081       0.266
MotherEd_3            0.1918      0.131      1.466      0.143      -0.065       0.449
MotherEd_4            0.1893      0.085      2.216      0.027       0.021       0.357
MotherEd_5            0.1747      0.098      1.779      0.076      -0.018       0.368
MotherEd_6           -0.0352      0.133     -0.266      0.791      -0.296       0.225
MotherEd_7            0.2870      0.114      2.523      0.012       0.063       0.511
Campus_1             -0.0111      0.156     -0.071      0.943      -0.318       0.296
PostulationType_1     0.2153      0.086      2.516      0.012       0.047       0.383
==============================================================================
Omnibus:                      114.631   Durbin-Watson:                   1.889
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               65.966
Skew:                           0.723   Prob(JB):                     4.74e-15
Kurtosis:                       2.075   Cond. No.                     2.88e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 2.88e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.130
Model:                            OLS   Adj. R-squared:                  0.093
Method:                 Least Squares   F-statistic:                     3.502
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           2.09e-07
Time:                        16:11:18   Log-Likelihood:                -303.41
No. Observations:                 537   AIC:                             652.8
Df Residuals:                     514   BIC:                             751.4
Df Model:                          22                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.1145      0.322      3.467      0.001       0.483       1.746
NEMScr               -0.0022      0.001     -1.913      0.056      -0.005    5.98e-05
Ranking               0.0013      0.001      1.409      0.159      -0.001       0.003
LangScr               0.0010      0.000      2.799      0.005       0.000       0.002
MathScr              -0.0021      0.001     -3.986      0.000      -0.003      -0.001
ScienScr          -6.406e-05      0.000     -0.156      0.876      -0.001       0.001
SchoolRegion_1        0.1036      0.068      1.518      0.130      -0.030       0.238
SchoolRegion_2        0.2794      0.091      3.073      0.002       0.101       0.458
SchoolRegion_3        0.2089      0.107      1.950      0.052      -0.002       0.419
EdTypeCode_1         -0.0630      0.110     -0.571      0.568      -0.280       0.154
EdTypeCode_2          0.1235      0.105      1.179      0.239      -0.082       0.329
SchoolType_1         -0.0099      0.111     -0.089      0.929      -0.229       0.209
SchoolType_2         -0.1413      0.096     -1.474      0.141      -0.330       0.047
SchoolType_3         -0.0675      0.091     -0.741      0.459      -0.246       0.111
SchoolType_4         -0.2114      0.281     -0.751      0.453      -0.764       0.341
MotherEd_1            0.0757      0.102      0.743      0.458      -0.124       0.276
MotherEd_2            0.0926      0.088      1.050      0.294      -0.081       0.266
MotherEd_3            0.1923      0.130      1.475      0.141      -0.064       0.448
MotherEd_4            0.1895      0.085      2.221      0.027       0.022       0.357
MotherEd_5            0.1749      0.098      1.783      0.075      -0.018       0.368
MotherEd_6           -0.0349      0.132     -0.264      0.792      -0.295       0.225
MotherEd_7            0.2871      0.114      2.526      0.012       0.064       0.510
PostulationType_1     0.2154      0.085      2.521      0.012       0.048       0.383
==============================================================================
Omnibus:                      114.718   Durbin-Watson:                   1.889
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               66.009
Skew:                           0.723   Prob(JB):                     4.64e-15
Kurtosis:                       2.074   Cond. No.                     2.22e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 2.22e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.130
Model:                            OLS   Adj. R-squared:                  0.095
Method:                 Least Squares   F-statistic:                     3.675
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           1.06e-07
Time:                        16:11:19   Log-Likelihood:                -303.41
No. Observations:                 537   AIC:                             650.8
Df Residuals:                     515   BIC:                             745.1
Df Model:                          21                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.1086      0.314      3.528      0.000       0.491       1.726
NEMScr               -0.0022      0.001     -1.919      0.056      -0.005     5.3e-05
Ranking               0.0013      0.001      1.412      0.159      -0.001       0.003
LangScr               0.0010      0.000      2.801      0.005       0.000       0.002
MathScr              -0.0021      0.001     -3.989      0.000      -0.003      -0.001
ScienScr          -6.291e-05      0.000     -0.153      0.878      -0.001       0.001
SchoolRegion_1        0.1043      0.068      1.539      0.124      -0.029       0.237
SchoolRegion_2        0.2796      0.091      3.079      0.002       0.101       0.458
SchoolRegion_3        0.2099      0.106      1.973      0.049       0.001       0.419
EdTypeCode_1         -0.0630      0.110     -0.572      0.568      -0.279       0.153
EdTypeCode_2          0.1228      0.104      1.177      0.240      -0.082       0.328
SchoolType_2         -0.1353      0.068     -1.998      0.046      -0.268      -0.002
SchoolType_3         -0.0614      0.060     -1.029      0.304      -0.179       0.056
SchoolType_4         -0.2052      0.272     -0.754      0.451      -0.740       0.330
MotherEd_1            0.0753      0.102      0.741      0.459      -0.124       0.275
MotherEd_2            0.0924      0.088      1.050      0.294      -0.081       0.265
MotherEd_3            0.1917      0.130      1.474      0.141      -0.064       0.447
MotherEd_4            0.1893      0.085      2.221      0.027       0.022       0.357
MotherEd_5            0.1748      0.098      1.784      0.075      -0.018       0.367
MotherEd_6           -0.0347      0.132     -0.262      0.793      -0.294       0.225
MotherEd_7            0.2868      0.113      2.527      0.012       0.064       0.510
PostulationType_1     0.2155      0.085      2.525      0.012       0.048       0.383
==============================================================================
Omnibus:                      114.797   Durbin-Watson:                   1.889
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               65.993
Skew:                           0.723   Prob(JB):                     4.68e-15
Kurtosis:                       2.074   Cond. No.                     2.15e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 2.15e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.130
Model:                            OLS   Adj. R-squared:                  0.097
Method:                 Least Squares   F-statistic:                     3.865
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           5.28e-08
Time:                        16:11:19   Log-Likelihood:                -303.43
No. Observations:                 537   AIC:                             648.9
Df Residuals:                     516   BIC:                             738.9
Df Model:                          20                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.1002      0.309      3.560      0.000       0.493       1.707
NEMScr               -0.0022      0.001     -1.918      0.056      -0.005    5.44e-05
Ranking               0.0013      0.001      1.408      0.160      -0.001       0.003
LangScr               0.0010      0.000      2.847      0.005       0.000       0.002
MathScr              -0.0021      0.000     -4.437      0.000      -0.003      -0.001
SchoolRegion_1        0.1048      0.068      1.549      0.122      -0.028       0.238
SchoolRegion_2        0.2801      0.091      3.089      0.002       0.102       0.458
SchoolRegion_3        0.2105      0.106      1.982      0.048       0.002       0.419
EdTypeCode_1         -0.0633      0.110     -0.575      0.565      -0.279       0.153
EdTypeCode_2          0.1239      0.104      1.192      0.234      -0.080       0.328
SchoolType_2         -0.1342      0.067     -1.995      0.047      -0.266      -0.002
SchoolType_3         -0.0601      0.059     -1.018      0.309      -0.176       0.056
SchoolType_4         -0.2027      0.271     -0.747      0.456      -0.736       0.331
MotherEd_1            0.0749      0.102      0.738      0.461      -0.125       0.274
MotherEd_2            0.0926      0.088      1.052      0.293      -0.080       0.265
MotherEd_3            0.1914      0.130      1.473      0.141      -0.064       0.447
MotherEd_4            0.1891      0.085      2.221      0.027       0.022       0.356
MotherEd_5            0.1746      0.098      1.784      0.075      -0.018       0.367
MotherEd_6           -0.0355      0.132     -0.269      0.788      -0.295       0.224
MotherEd_7            0.2877      0.113      2.541      0.011       0.065       0.510
PostulationType_1     0.2156      0.085      2.528      0.012       0.048       0.383
==============================================================================
Omnibus:                      114.765   Durbin-Watson:                   1.890
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               65.943
Skew:                           0.723   Prob(JB):                     4.79e-15
Kurtosis:                       2.074   Cond. No.                     1.90e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.9e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.130
Model:                            OLS   Adj. R-squared:                  0.098
Method:                 Least Squares   F-statistic:                     4.072
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           2.61e-08
Time:                        16:11:19   Log-Likelihood:                -303.46
No. Observations:                 537   AIC:                             646.9
Df Residuals:                     517   BIC:                             732.6
Df Model:                          19                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0854      0.304      3.572      0.000       0.488       1.682
NEMScr               -0.0022      0.001     -1.914      0.056      -0.004    5.82e-05
Ranking               0.0013      0.001      1.408      0.160      -0.001       0.003
LangScr               0.0010      0.000      2.843      0.005       0.000       0.002
MathScr              -0.0021      0.000     -4.440      0.000      -0.003      -0.001
SchoolRegion_1        0.1057      0.067      1.566      0.118      -0.027       0.238
SchoolRegion_2        0.2808      0.091      3.101      0.002       0.103       0.459
SchoolRegion_3        0.2109      0.106      1.988      0.047       0.003       0.419
EdTypeCode_1         -0.0642      0.110     -0.585      0.559      -0.280       0.152
EdTypeCode_2          0.1241      0.104      1.195      0.233      -0.080       0.328
SchoolType_2         -0.1347      0.067     -2.006      0.045      -0.267      -0.003
SchoolType_3         -0.0601      0.059     -1.020      0.308      -0.176       0.056
SchoolType_4         -0.2024      0.271     -0.746      0.456      -0.735       0.330
MotherEd_1            0.0876      0.090      0.974      0.330      -0.089       0.264
MotherEd_2            0.1050      0.075      1.401      0.162      -0.042       0.252
MotherEd_3            0.2040      0.121      1.683      0.093      -0.034       0.442
MotherEd_4            0.2015      0.072      2.816      0.005       0.061       0.342
MotherEd_5            0.1872      0.086      2.181      0.030       0.019       0.356
MotherEd_7            0.2999      0.104      2.893      0.004       0.096       0.504
PostulationType_1     0.2154      0.085      2.528      0.012       0.048       0.383
==============================================================================
Omnibus:                      114.283   Durbin-Watson:                   1.890
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               66.024
Skew:                           0.724   Prob(JB):                     4.60e-15
Kurtosis:                       2.076   Cond. No.                     1.83e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.83e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(6)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.130
Model:                            OLS   Adj. R-squared:                  0.099
Method:                 Least Squares   F-statistic:                     4.285
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           1.40e-08
Time:                        16:11:19   Log-Likelihood:                -303.64
No. Observations:                 537   AIC:                             645.3
Df Residuals:                     518   BIC:                             726.7
Df Model:                          18                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0742      0.303      3.544      0.000       0.479       1.670
NEMScr               -0.0021      0.001     -1.845      0.066      -0.004       0.000
Ranking               0.0012      0.001      1.323      0.186      -0.001       0.003
LangScr               0.0010      0.000      2.835      0.005       0.000       0.002
MathScr              -0.0021      0.000     -4.417      0.000      -0.003      -0.001
SchoolRegion_1        0.1037      0.067      1.539      0.124      -0.029       0.236
SchoolRegion_2        0.2784      0.090      3.079      0.002       0.101       0.456
SchoolRegion_3        0.2100      0.106      1.980      0.048       0.002       0.418
EdTypeCode_2          0.1308      0.103      1.268      0.206      -0.072       0.333
SchoolType_2         -0.1357      0.067     -2.022      0.044      -0.267      -0.004
SchoolType_3         -0.0629      0.059     -1.071      0.285      -0.178       0.053
SchoolType_4         -0.2412      0.263     -0.918      0.359      -0.758       0.275
MotherEd_1            0.0886      0.090      0.986      0.325      -0.088       0.265
MotherEd_2            0.1050      0.075      1.402      0.161      -0.042       0.252
MotherEd_3            0.2002      0.121      1.656      0.098      -0.037       0.438
MotherEd_4            0.2030      0.071      2.841      0.005       0.063       0.343
MotherEd_5            0.1899      0.086      2.216      0.027       0.022       0.358
MotherEd_7            0.3005      0.104      2.900      0.004       0.097       0.504
PostulationType_1     0.2169      0.085      2.548      0.011       0.050       0.384
==============================================================================
Omnibus:                      114.772   Durbin-Watson:                   1.882
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               65.551
Skew:                           0.720   Prob(JB):                     5.83e-15
Kurtosis:                       2.073   Cond. No.                     1.83e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.83e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(7)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.128
Model:                            OLS   Adj. R-squared:                  0.100
Method:                 Least Squares   F-statistic:                     4.489
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           8.96e-09
Time:                        16:11:19   Log-Likelihood:                -304.08
No. Observations:                 537   AIC:                             644.2
Df Residuals:                     519   BIC:                             721.3
Df Model:                          17                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0676      0.303      3.524      0.000       0.472       1.663
NEMScr               -0.0020      0.001     -1.751      0.081      -0.004       0.000
Ranking               0.0011      0.001      1.213      0.226      -0.001       0.003
LangScr               0.0010      0.000      2.801      0.005       0.000       0.002
MathScr              -0.0021      0.000     -4.399      0.000      -0.003      -0.001
SchoolRegion_1        0.1020      0.067      1.515      0.130      -0.030       0.234
SchoolRegion_2        0.2778      0.090      3.073      0.002       0.100       0.455
SchoolRegion_3        0.2001      0.105      1.897      0.058      -0.007       0.407
EdTypeCode_2          0.1223      0.103      1.190      0.234      -0.080       0.324
SchoolType_2         -0.1260      0.066     -1.901      0.058      -0.256       0.004
SchoolType_3         -0.0528      0.058     -0.914      0.361      -0.166       0.061
MotherEd_1            0.0858      0.090      0.956      0.340      -0.091       0.262
MotherEd_2            0.1020      0.075      1.364      0.173      -0.045       0.249
MotherEd_3            0.1987      0.121      1.644      0.101      -0.039       0.436
MotherEd_4            0.2006      0.071      2.810      0.005       0.060       0.341
MotherEd_5            0.1881      0.086      2.197      0.028       0.020       0.356
MotherEd_7            0.3012      0.104      2.908      0.004       0.098       0.505
PostulationType_1     0.2154      0.085      2.531      0.012       0.048       0.383
==============================================================================
Omnibus:                      116.832   Durbin-Watson:                   1.888
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               65.414
Skew:                           0.716   Prob(JB):                     6.24e-15
Kurtosis:                       2.067   Cond. No.                     1.83e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.83e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(8)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.127
Model:                            OLS   Adj. R-squared:                  0.100
Method:                 Least Squares   F-statistic:                     4.718
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           5.60e-09
Time:                        16:11:19   Log-Likelihood:                -304.51
No. Observations:                 537   AIC:                             643.0
Df Residuals:                     520   BIC:                             715.9
Df Model:                          16                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0139      0.297      3.412      0.001       0.430       1.598
NEMScr               -0.0020      0.001     -1.742      0.082      -0.004       0.000
Ranking               0.0011      0.001      1.195      0.233      -0.001       0.003
LangScr               0.0010      0.000      2.853      0.005       0.000       0.002
MathScr              -0.0021      0.000     -4.388      0.000      -0.003      -0.001
SchoolRegion_1        0.0999      0.067      1.485      0.138      -0.032       0.232
SchoolRegion_2        0.2784      0.090      3.080      0.002       0.101       0.456
SchoolRegion_3        0.1988      0.105      1.886      0.060      -0.008       0.406
EdTypeCode_2          0.1283      0.103      1.252      0.211      -0.073       0.330
SchoolType_2         -0.0819      0.045     -1.802      0.072      -0.171       0.007
MotherEd_1            0.0795      0.090      0.889      0.375      -0.096       0.255
MotherEd_2            0.1015      0.075      1.358      0.175      -0.045       0.248
MotherEd_3            0.1938      0.121      1.605      0.109      -0.043       0.431
MotherEd_4            0.2003      0.071      2.806      0.005       0.060       0.340
MotherEd_5            0.1830      0.085      2.142      0.033       0.015       0.351
MotherEd_7            0.3057      0.103      2.955      0.003       0.102       0.509
PostulationType_1     0.2197      0.085      2.587      0.010       0.053       0.387
==============================================================================
Omnibus:                      123.139   Durbin-Watson:                   1.891
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               65.590
Skew:                           0.712   Prob(JB):                     5.72e-15
Kurtosis:                       2.050   Cond. No.                     1.79e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.79e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(9)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.125
Model:                            OLS   Adj. R-squared:                  0.100
Method:                 Least Squares   F-statistic:                     4.982
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           3.35e-09
Time:                        16:11:19   Log-Likelihood:                -304.92
No. Observations:                 537   AIC:                             641.8
Df Residuals:                     521   BIC:                             710.4
Df Model:                          15                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0528      0.294      3.583      0.000       0.476       1.630
NEMScr               -0.0020      0.001     -1.734      0.084      -0.004       0.000
Ranking               0.0011      0.001      1.185      0.237      -0.001       0.003
LangScr               0.0010      0.000      2.810      0.005       0.000       0.002
MathScr              -0.0021      0.000     -4.363      0.000      -0.003      -0.001
SchoolRegion_1        0.1005      0.067      1.494      0.136      -0.032       0.233
SchoolRegion_2        0.2776      0.090      3.072      0.002       0.100       0.455
SchoolRegion_3        0.1969      0.105      1.868      0.062      -0.010       0.404
EdTypeCode_2          0.1244      0.102      1.215      0.225      -0.077       0.326
SchoolType_2         -0.0772      0.045     -1.710      0.088      -0.166       0.012
MotherEd_2            0.0625      0.060      1.033      0.302      -0.056       0.181
MotherEd_3            0.1556      0.113      1.379      0.168      -0.066       0.377
MotherEd_4            0.1607      0.056      2.880      0.004       0.051       0.270
MotherEd_5            0.1424      0.072      1.973      0.049       0.001       0.284
MotherEd_7            0.2661      0.093      2.851      0.005       0.083       0.450
PostulationType_1     0.2199      0.085      2.590      0.010       0.053       0.387
==============================================================================
Omnibus:                      125.837   Durbin-Watson:                   1.892
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               66.034
Skew:                           0.714   Prob(JB):                     4.58e-15
Kurtosis:                       2.044   Cond. No.                     1.76e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.76e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(10)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.124
Model:                            OLS   Adj. R-squared:                  0.100
Method:                 Least Squares   F-statistic:                     5.261
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           2.20e-09
Time:                        16:11:19   Log-Likelihood:                -305.47
No. Observations:                 537   AIC:                             640.9
Df Residuals:                     522   BIC:                             705.2
Df Model:                          14                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0790      0.293      3.686      0.000       0.504       1.654
NEMScr               -0.0019      0.001     -1.682      0.093      -0.004       0.000
Ranking               0.0011      0.001      1.152      0.250      -0.001       0.003
LangScr               0.0010      0.000      2.768      0.006       0.000       0.002
MathScr              -0.0021      0.000     -4.362      0.000      -0.003      -0.001
SchoolRegion_1        0.1027      0.067      1.528      0.127      -0.029       0.235
SchoolRegion_2        0.2858      0.090      3.176      0.002       0.109       0.463
SchoolRegion_3        0.2010      0.105      1.909      0.057      -0.006       0.408
EdTypeCode_2          0.1213      0.102      1.185      0.237      -0.080       0.322
SchoolType_2         -0.0824      0.045     -1.836      0.067      -0.171       0.006
MotherEd_3            0.1204      0.108      1.119      0.264      -0.091       0.332
MotherEd_4            0.1263      0.045      2.822      0.005       0.038       0.214
MotherEd_5            0.1094      0.065      1.690      0.092      -0.018       0.237
MotherEd_7            0.2306      0.087      2.657      0.008       0.060       0.401
PostulationType_1     0.2190      0.085      2.579      0.010       0.052       0.386
==============================================================================
Omnibus:                      128.802   Durbin-Watson:                   1.892
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               66.515
Skew:                           0.715   Prob(JB):                     3.60e-15
Kurtosis:                       2.037   Cond. No.                     1.76e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.76e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(11)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.122
Model:                            OLS   Adj. R-squared:                  0.100
Method:                 Least Squares   F-statistic:                     5.567
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           1.52e-09
Time:                        16:11:19   Log-Likelihood:                -306.11
No. Observations:                 537   AIC:                             640.2
Df Residuals:                     523   BIC:                             700.2
Df Model:                          13                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.1204      0.290      3.857      0.000       0.550       1.691
NEMScr               -0.0019      0.001     -1.646      0.100      -0.004       0.000
Ranking               0.0010      0.001      1.109      0.268      -0.001       0.003
LangScr               0.0010      0.000      2.795      0.005       0.000       0.002
MathScr              -0.0022      0.000     -4.523      0.000      -0.003      -0.001
SchoolRegion_1        0.1034      0.067      1.538      0.125      -0.029       0.236
SchoolRegion_2        0.2859      0.090      3.176      0.002       0.109       0.463
SchoolRegion_3        0.1989      0.105      1.889      0.059      -0.008       0.406
EdTypeCode_2          0.1237      0.102      1.209      0.227      -0.077       0.325
SchoolType_2         -0.0854      0.045     -1.906      0.057      -0.173       0.003
MotherEd_4            0.1178      0.044      2.669      0.008       0.031       0.204
MotherEd_5            0.1019      0.064      1.582      0.114      -0.025       0.228
MotherEd_7            0.2217      0.086      2.565      0.011       0.052       0.392
PostulationType_1     0.2189      0.085      2.577      0.010       0.052       0.386
==============================================================================
Omnibus:                      130.086   Durbin-Watson:                   1.886
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               67.197
Skew:                           0.720   Prob(JB):                     2.56e-15
Kurtosis:                       2.035   Cond. No.                     1.74e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.74e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(12)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.119
Model:                            OLS   Adj. R-squared:                  0.099
Method:                 Least Squares   F-statistic:                     5.926
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           1.02e-09
Time:                        16:11:19   Log-Likelihood:                -306.74
No. Observations:                 537   AIC:                             639.5
Df Residuals:                     524   BIC:                             695.2
Df Model:                          12                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.0990      0.290      3.791      0.000       0.530       1.669
NEMScr               -0.0006      0.000     -2.445      0.015      -0.001      -0.000
LangScr               0.0010      0.000      2.863      0.004       0.000       0.002
MathScr              -0.0023      0.000     -5.144      0.000      -0.003      -0.001
SchoolRegion_1        0.1112      0.067      1.661      0.097      -0.020       0.243
SchoolRegion_2        0.2855      0.090      3.171      0.002       0.109       0.462
SchoolRegion_3        0.2037      0.105      1.935      0.054      -0.003       0.410
EdTypeCode_2          0.1464      0.100      1.460      0.145      -0.051       0.343
SchoolType_2         -0.0877      0.045     -1.959      0.051      -0.176       0.000
MotherEd_4            0.1213      0.044      2.757      0.006       0.035       0.208
MotherEd_5            0.1042      0.064      1.617      0.106      -0.022       0.231
MotherEd_7            0.2208      0.086      2.554      0.011       0.051       0.391
PostulationType_1     0.2234      0.085      2.632      0.009       0.057       0.390
==============================================================================
Omnibus:                      133.703   Durbin-Watson:                   1.883
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               67.886
Skew:                           0.723   Prob(JB):                     1.81e-15
Kurtosis:                       2.028   Cond. No.                     1.51e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.51e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(13)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.116
Model:                            OLS   Adj. R-squared:                  0.097
Method:                 Least Squares   F-statistic:                     6.257
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           9.78e-10
Time:                        16:11:19   Log-Likelihood:                -307.83
No. Observations:                 537   AIC:                             639.7
Df Residuals:                     525   BIC:                             691.1
Df Model:                          11                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.1538      0.288      4.009      0.000       0.588       1.719
NEMScr               -0.0007      0.000     -2.487      0.013      -0.001      -0.000
LangScr               0.0010      0.000      2.806      0.005       0.000       0.002
MathScr              -0.0024      0.000     -5.222      0.000      -0.003      -0.001
SchoolRegion_1        0.1049      0.067      1.570      0.117      -0.026       0.236
SchoolRegion_2        0.2763      0.090      3.073      0.002       0.100       0.453
SchoolRegion_3        0.1919      0.105      1.826      0.068      -0.015       0.398
SchoolType_2         -0.0886      0.045     -1.976      0.049      -0.177      -0.001
MotherEd_4            0.1242      0.044      2.823      0.005       0.038       0.211
MotherEd_5            0.1013      0.064      1.572      0.117      -0.025       0.228
MotherEd_7            0.2272      0.086      2.628      0.009       0.057       0.397
PostulationType_1     0.2159      0.085      2.546      0.011       0.049       0.383
==============================================================================
Omnibus:                      142.955   Durbin-Watson:                   1.896
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               68.561
Skew:                           0.721   Prob(JB):                     1.29e-15
Kurtosis:                       2.008   Cond. No.                     1.50e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.5e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(14)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.112
Model:                            OLS   Adj. R-squared:                  0.095
Method:                 Least Squares   F-statistic:                     6.618
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           1.06e-09
Time:                        16:11:19   Log-Likelihood:                -309.09
No. Observations:                 537   AIC:                             640.2
Df Residuals:                     526   BIC:                             687.3
Df Model:                          10                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.2135      0.286      4.248      0.000       0.652       1.775
NEMScr               -0.0007      0.000     -2.589      0.010      -0.001      -0.000
LangScr               0.0010      0.000      2.792      0.005       0.000       0.002
MathScr              -0.0023      0.000     -5.091      0.000      -0.003      -0.001
SchoolRegion_2        0.1839      0.068      2.702      0.007       0.050       0.318
SchoolRegion_3        0.0982      0.087      1.134      0.257      -0.072       0.268
SchoolType_2         -0.0877      0.045     -1.955      0.051      -0.176       0.000
MotherEd_4            0.1271      0.044      2.886      0.004       0.041       0.214
MotherEd_5            0.1023      0.065      1.585      0.113      -0.024       0.229
MotherEd_7            0.2350      0.086      2.719      0.007       0.065       0.405
PostulationType_1     0.2256      0.085      2.664      0.008       0.059       0.392
==============================================================================
Omnibus:                      139.289   Durbin-Watson:                   1.897
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               70.266
Skew:                           0.738   Prob(JB):                     5.52e-16
Kurtosis:                       2.019   Cond. No.                     1.48e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.48e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(15)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.110
Model:                            OLS   Adj. R-squared:                  0.094
Method:                 Least Squares   F-statistic:                     7.207
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           6.67e-10
Time:                        16:11:19   Log-Likelihood:                -309.74
No. Observations:                 537   AIC:                             639.5
Df Residuals:                     527   BIC:                             682.3
Df Model:                           9                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.2225      0.286      4.280      0.000       0.661       1.784
NEMScr               -0.0007      0.000     -2.544      0.011      -0.001      -0.000
LangScr               0.0010      0.000      2.782      0.006       0.000       0.002
MathScr              -0.0023      0.000     -5.132      0.000      -0.003      -0.001
SchoolRegion_2        0.1783      0.068      2.626      0.009       0.045       0.312
SchoolType_2         -0.0839      0.045     -1.875      0.061      -0.172       0.004
MotherEd_4            0.1293      0.044      2.939      0.003       0.043       0.216
MotherEd_5            0.1042      0.065      1.615      0.107      -0.023       0.231
MotherEd_7            0.2342      0.086      2.709      0.007       0.064       0.404
PostulationType_1     0.2258      0.085      2.666      0.008       0.059       0.392
==============================================================================
Omnibus:                      145.490   Durbin-Watson:                   1.891
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               70.704
Skew:                           0.737   Prob(JB):                     4.43e-16
Kurtosis:                       2.006   Cond. No.                     1.48e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.48e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'

$ws = $wb.Worksheets.Item(16)
$ws.Range("B2").Value = '                            OLS Regression Results                            
==============================================================================
Dep. Variable:               Desertor   R-squared:                       0.101
Model:                            OLS   Adj. R-squared:                  0.089
Method:                 Least Squares   F-statistic:                     8.487
Date:                Sun, 29 Dec 2019   Prob (F-statistic):           7.28e-10
Time:                        16:11:19   Log-Likelihood:                -312.33
No. Observations:                 537   AIC:                             640.7
Df Residuals:                     529   BIC:                             674.9
Df Model:                           7                                         
Covariance Type:            nonrobust                                         
=====================================================================================
                        coef    std err          t      P>|t|      [0.025      0.975]
-------------------------------------------------------------------------------------
Ones                  1.2218      0.286      4.278      0.000       0.661       1.783
NEMScr               -0.0006      0.000     -2.336      0.020      -0.001   -9.74e-05
LangScr               0.0010      0.000      2.984      0.003       0.000       0.002
MathScr              -0.0024      0.000     -5.458      0.000      -0.003      -0.002
SchoolRegion_2        0.1739      0.068      2.556      0.011       0.040       0.308
MotherEd_4            0.1020      0.042      2.448      0.015       0.020       0.184
MotherEd_7            0.2104      0.086      2.456      0.014       0.042       0.379
PostulationType_1     0.2341      0.085      2.760      0.006       0.067       0.401
==============================================================================
Omnibus:                      147.447   Durbin-Watson:                   1.897
Prob(Omnibus):                  0.000   Jarque-Bera (JB):               73.833
Skew:                           0.760   Prob(JB):                     9.28e-17
Kurtosis:                       2.007   Cond. No.                     1.48e+04
==============================================================================

Warnings:
[1] Standard Errors assume that the covariance matrix of the errors is correctly specified.
[2] The condition number is large, 1.48e+04. This might indicate that there are
strong multicollinearity or other numerical problems.'
